$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "26.655.40"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "1.630.94"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0839"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("D12").Value = "1.859.75"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").Value = "1.612.68"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").Value = "26.633.23"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.09%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("E24").Value = "  +3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0520"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.05%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "1.167.52"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.808"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.792"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "1.769.19"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0509"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("E49").Value = "  +4.46%  "
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  +0.16%  "
